# Insert a new weekly price record (row 37) into the Haba price series.
# This shifts all existing rows from 37..59 down to 38..60, extending the
# sheet's dimension from A1:R59 to A1:R60, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 37.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record's data.
$ws.Range("A37").Value = 4
$ws.Range("B37").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C37").Value = "Los Lagos"
$ws.Range("D37").Value = 44518
$ws.Range("E37").Value = 10
$ws.Range("F37").Value = 100112026
$ws.Range("G37").Value = "Haba"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 80
$ws.Range("K37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = 10000
$ws.Range("N37").Value = "`$/saco 25 kilos"
$ws.Range("O37").Value = "Región del Maule"
$ws.Range("P37").Value = 400
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
